# Update the "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the latest generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# Sheet "展览" -> F column updates (by row number)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 309
$ws1.Range("F9").Value = 24
$ws1.Range("F13").Value = 2956
$ws1.Range("F19").Value = 555
$ws1.Range("F23").Value = 105
$ws1.Range("F24").Value = 53
$ws1.Range("F27").Value = 2291
$ws1.Range("F28").Value = 4837

# Sheet "全部类型" -> F column updates (by row number)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 309
$ws4.Range("F9").Value = 24
$ws4.Range("F13").Value = 2956
$ws4.Range("F20").Value = 555
$ws4.Range("F24").Value = 105
$ws4.Range("F25").Value = 53
$ws4.Range("F28").Value = 2291
$ws4.Range("F29").Value = 4837
